$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new summary-percentage cells in row 13
$ws.Range("E13").Formula = "=SUM(E3:E12)/COUNT(E3:E12)"
$ws.Range("F13").Formula = "=SUM(F3:F12)/COUNT(F3:F12)"

# Match the resulting view state (ruler shown, selection moved)
$excel.ActiveWindow.DisplayRuler = $true
$ws.Range("C16").Select()
